{"js": "// SAE Word template fix: replace the three hard-coded, numbered \"Exemple N : \"\n// list paragraphs with a single plain paragraph containing the ${exemples}\n// merge-field placeholder (plus a trailing space), matching the move away\n// from static examples towards a template token.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate every paragraph that starts the bulleted \"Exemple N : \" list items.\nconst exampleIndexes = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Exemple \") === 0) {\n    exampleIndexes.push(i);\n  }\n}\n\nif (exampleIndexes.length > 0) {\n  // Anchor on the paragraph right before the first \"Exemple\" item (the\n  // \"Liste d'exemples de SA\u00c9 :\" paragraph) and insert a brand-new plain\n  // paragraph after it, so it does NOT inherit the \"Paragraphedeliste\"\n  // list style / numbering that the \"Exemple\" paragraphs use.\n  const anchor = paragraphs.items[exampleIndexes[0] - 1];\n  const newParagraph = anchor.insertParagraph(\"${exemples}\", Word.InsertLocation.after);\n  newParagraph.insertText(\" \", Word.InsertLocation.end);\n\n  // Delete the original \"Exemple N : \" list paragraphs, from last to first\n  // so earlier indexes stay valid while deleting.\n  for (let k = exampleIndexes.length - 1; k >= 0; k--) {\n    paragraphs.items[exampleIndexes[k]].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# SAE Word template fix: replace the three hard-coded, numbered \"Exemple N : \"\n# list paragraphs with a single plain paragraph containing the ${exemples}\n# merge-field placeholder (plus a trailing space), matching the move away\n# from static examples towards a template token.\n\n$d = $word.ActiveDocument\n\n# Locate every paragraph that starts the bulleted \"Exemple N : \" list items.\n$exampleIdxs = @()\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.StartsWith(\"Exemple \")) {\n        $exampleIdxs += $i\n    }\n}\n\nif ($exampleIdxs.Count -gt 0) {\n    $firstIdx = $exampleIdxs[0]\n\n    # Insert a brand-new plain paragraph right before the first \"Exemple\" item,\n    # anchored on the preceding paragraph so the new paragraph does NOT inherit\n    # the \"Paragraphedeliste\" list style / numbering.\n    $anchor = $d.Paragraphs($firstIdx - 1)\n    $anchor.Range.InsertParagraphAfter()\n\n    # The freshly inserted (empty) paragraph is now at $firstIdx.\n    $newPara = $d.Paragraphs($firstIdx)\n    $newPara.Range.Text = '${exemples} '\n\n    # Delete the original \"Exemple N : \" list paragraphs (now shifted down by\n    # one because of the insertion above), from last to first so indices stay\n    # valid while deleting.\n    for ($k = $exampleIdxs.Count - 1; $k -ge 0; $k--) {\n        $d.Paragraphs($exampleIdxs[$k] + 1).Range.Delete()\n    }\n}\n"}
